# Update the "Pure Luck" column (F) values on Sheet1 to a freshly
# re-rolled set of random numbers (new input mode for luck factor).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$luckValues = @{
    2  = 0.9816078898292088
    3  = 0.3093259448199278
    4  = 0.528011429751901
    5  = 0.2211431802606283
    6  = 0.2025396008973027
    7  = 0.3353222317810357
    8  = 0.480881821685384
    9  = 0.2456892014419062
    10 = 0.6399064000863915
    11 = 0.3897293990344669
    12 = 0.4177155948873165
    13 = 0.3291603309168202
    14 = 0.2592286777889876
    15 = 0.2509383021706575
    16 = 0.2636604920511952
    17 = 0.7188558658383181
    18 = 0.2305833645066783
}

foreach ($row in $luckValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $luckValues[$row]
}
